# "f string notes released"
# - bump the auto date field (datetimeFigureOut) on the slide master and every
#   slide layout from 2024-09-10 to 2024-09-12 (PowerPoint re-stamps these
#   fields on save)
# - tighten the printed precision of the three pi format-spec textboxes on
#   slide 10 (f"{pi:>10.7}" / f"{pi:^10.7}" / f"{pi:<10.7}") to the "6f" form,
#   and the two on slide 9 (f"{pi:10.9}" / f"{pi:10.7}") to the "8f"/"6f" form
# - these Courier-New textboxes are "resize shape to fit text" (spAutoFit),
#   so when the text changes width/height PowerPoint relayouts the box; we
#   reproduce the exact resulting geometry here

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder refresh: master + all 11 slide layouts
# ---------------------------------------------------------------------
$m = $p.SlideMaster

function Set-DatePlaceholderText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "2024-09-12"
        }
    }
}

Set-DatePlaceholderText $m.Shapes

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $cl = $m.CustomLayouts.Item($li)
    Set-DatePlaceholderText $cl.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 10 (index 10): three pi format textboxes, precision 7 -> 6f
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)

# TextBox 3 : f"{pi:>10.7}" -> f"{pi:>10.6f}"  (off x nudges by -1 EMU, ext grows)
$tb3 = $s10.Shapes.Item(1)
$tb3.TextFrame.TextRange.Text = 'f"{pi:>10.6f}"'
$tb3.Left = 627.87302
$tb3.Width = 311.22371

# TextBox 5 : f"{pi:^10.7}" -> f"{pi:^10.6f}"  (only the width grows)
$tb5 = $s10.Shapes.Item(3)
$tb5.TextFrame.TextRange.Text = 'f"{pi:^10.6f}"'
$tb5.Width = 311.22371

# TextBox 8 : f"{pi:<10.7}" -> f"{pi:<10.6f}" (box geometry unchanged)
$tb8 = $s10.Shapes.Item(5)
$tb8.TextFrame.TextRange.Text = 'f"{pi:<10.6f}"'

# ---------------------------------------------------------------------
# 3) Slide 9 (index 9): two pi format textboxes
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)

# TextBox 6 : f"{pi:10.9}" -> f"{pi:10.8f}" (text now wraps to two lines, box grows taller)
$tb6 = $s9.Shapes.Item(2)
$tb6.TextFrame.TextRange.Text = 'f"{pi:10.8f}"'
$tb6.Height = 84.82032

# TextBox 7 : f"{pi:10.7}" -> f"{pi:10.6f}" (off x nudges by -1 EMU, width grows)
$tb7 = $s9.Shapes.Item(3)
$tb7.TextFrame.TextRange.Text = 'f"{pi:10.6f}"'
$tb7.Left = 204.3892
$tb7.Width = 275.6108
